$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header info (account holder name / card number) ---
$ws.Range("C2").Value = "Hartmut"

# The card number is all digits, so a plain .Value assignment gets
# auto-typed as a number by Excel. Force it to stay text (matching the
# original cell, which stored it as a string) by formatting as Text
# first, then normalize the number format back via a neighbouring
# same-styled cell so the style index itself is unaffected.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance date ---
$ws.Range("D5").Value = "KONTOSTAND AM 12.09.2023"

# --- Transaction rows ---
$ws.Range("B6").Value = "13.09."
$ws.Range("C6").Value = "14.09."
$ws.Range("D6").Value = "EBAY MKTPLC EU SSIYLG"
$ws.Range("E6").Value = "186,54-"

$ws.Range("B7").Value = "14.09."
$ws.Range("C7").Value = "15.09."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "84,15-"

$ws.Range("B8").Value = "18.09."
$ws.Range("C8").Value = "19.09."
$ws.Range("D8").Value = "PAYPAL ZXTGMG"
$ws.Range("E8").Value = "92,03-"

$ws.Range("B9").Value = "20.09."
$ws.Range("C9").Value = "21.09."
$ws.Range("D9").Value = "RECHNUNG VODAFONE GMBH 97417508"
$ws.Range("E9").Value = "40,48-"

# Row 10 transaction removed entirely - clear values & match the blank
# spacer-row formatting already used by row 11 (E column there keeps its
# "right" alignment but also picks up vertical-center + wrap, same as the
# other blank spacer row).
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152
$ws.Range("E10").VerticalAlignment = -4108
$ws.Range("E10").WrapText = $true

# --- Closing balance / next statement date ---
$ws.Range("D12").Value = "KONTOSTAND AM 25.09.2023"
$ws.Range("E12").Value = "403,20-"
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 03.10.2023"
